$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force text formatting so numeric-looking strings (CNPJ/CPF)
# and the timestamp string are written as text, not auto-converted to
# numbers/dates by Excel.
$ws.Range("A2:E2").NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"

$ws.Range("A2").Value = "43213049000135"
$ws.Range("B2").Value = "RHAIANY SOARES SILVA"
$ws.Range("C2").Value = "14707320767"
$ws.Range("D2").Value = "TUR-PRO-2025/01056"
$ws.Range("E2").Value = "https://acesso.processo.rio/sigaex/public/app/transparencia/processo?n=TUR-PRO-2025/01056"
$ws.Range("H2").Value = "2025-12-29T16:57:53.530794"

# Restore the cells' style to the workbook default so no visible
# formatting change is introduced (matches original unstyled cells).
$ws.Range("A2:E2").Style = "Normal"
$ws.Range("H2").Style = "Normal"
